$wb = $excel.ActiveWorkbook

$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
$wsWing = $wb.Worksheets.Item("WING")
$wsHTail = $wb.Worksheets.Item("HORIZONTAL TAIL")
$wsVTail = $wb.Worksheets.Item("VERTICAL TAIL")

# --- FUSELAGE ---
$wsFuselage.Range("C8").Value = 4156.0
$wsFuselage.Range("D8").Value = 56.063419103720584
$wsFuselage.Range("C9").Value = 3698.0
$wsFuselage.Range("D9").Value = 38.864899866592566
$wsFuselage.Range("C10").Value = 5004.0
$wsFuselage.Range("D10").Value = 87.9069656388397
$wsFuselage.Range("C11").Value = 3106.0
$wsFuselage.Range("D11").Value = 16.634499455282995
$wsFuselage.Range("C12").Value = 3117.0
$wsFuselage.Range("D12").Value = 17.047564327790436
$wsFuselage.Range("C13").Value = 2491.0
$wsFuselage.Range("D13").Value = -6.459582053087591
$wsFuselage.Range("C14").Value = 4094.0
$wsFuselage.Range("D14").Value = 53.73523527686046

# --- WING ---
$wsWing.Range("A9").Value = "TORENBEEK_1982"
$wsWing.Range("C9").Value = 2426.0
$wsWing.Range("D9").Value = -14.056999939703871
$wsWing.Range("A11").Value = "KROO"
$wsWing.Range("C11").Value = 2311.0
$wsWing.Range("D11").Value = -18.130967378670917

# --- HORIZONTAL TAIL ---
$wsHTail.Range("A8").Value = "ROSKAM"
$wsHTail.Range("C8").Value = 251.0
$wsHTail.Range("D8").Value = -18.04009249061837
$wsHTail.Range("A9").Value = "JENKINSON"
$wsHTail.Range("C9").Value = 293.0
$wsHTail.Range("D9").Value = -4.32568565637921
$wsHTail.Range("A10").Value = "HOWE"
$wsHTail.Range("C10").Value = 238.0
$wsHTail.Range("D10").Value = -22.28502793931144
$wsHTail.Range("A11").Value = "NICOLAI_2013"
$wsHTail.Range("C11").Value = 122.0
$wsHTail.Range("D11").Value = -60.162913481495785
$wsHTail.Range("A12").Value = "RAYMER"
$wsHTail.Range("C12").Value = 142.0
$wsHTail.Range("D12").Value = -53.632243560429515
$wsHTail.Range("A13").Value = "SADRAEY"
$wsHTail.Range("C13").Value = 271.0
$wsHTail.Range("D13").Value = -11.509422569552102
$wsHTail.Range("C14").Value = 276.0
$wsHTail.Range("D14").Value = -9.876755089285535

# --- VERTICAL TAIL ---
$wsVTail.Range("A8").Value = "ROSKAM"
$wsVTail.Range("C8").Value = 278.0
$wsVTail.Range("D8").Value = -9.22368809717891
$wsVTail.Range("A9").Value = "JENKINSON"
$wsVTail.Range("C9").Value = 330.0
$wsVTail.Range("D9").Value = 7.756053697593382
$wsVTail.Range("A10").Value = "HOWE"
$wsVTail.Range("C10").Value = 512.0
$wsVTail.Range("D10").Value = 67.1851499792964
$wsVTail.Range("C11").Value = 88.0
$wsVTail.Range("D11").Value = -71.26505234730843
$wsVTail.Range("C12").Value = 419.0
$wsVTail.Range("D12").Value = 36.81753484633826
$wsVTail.Range("C13").Value = 391.0
$wsVTail.Range("D13").Value = 27.67459695684549

